$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Caller section:
#    - prefix the caller-name paragraph with a bold "Caller " label
#      (rendered as "Caller" (bold) + " " + name, three separate runs)
#    - split out ", CALLERREGION GOES HERE" into its own paragraph,
#      dropping the leading ", "
# ---------------------------------------------------------------------
$rngCaller = $d.Content
$rngCaller.Find.Execute("CALLERFIRSTNAME GOES HERE CALLERLASTNAME GOES HERE") | Out-Null
$callerStart = $rngCaller.Start

$insertPoint = $d.Range($callerStart, $callerStart)
$insertPoint.InsertBefore("Caller ")

$boldRange = $d.Range($callerStart, $callerStart + 6)
$boldRange.Bold = 1

$d.Content.Find.Execute(", CALLERREGION GOES HERE", $true, $false, $false, $false, $false, $true, 1, $false, "^pCALLERREGION GOES HERE", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Patient section:
#    - insert a new paragraph "Patient " (bold) + patient name before
#      the (bookmarked) paragraph that used to hold that name
#    - insert "PATIENTAGENCYID GOES HERE" and "PATIENTGENDER GOES HERE"
#      as their own paragraphs
#    - the bookmarked paragraph keeps the bookmark but ends up holding
#      only "PATIENTAGE GOES HERE"
#    - the old combined Agency/Gender/Age paragraph is emptied out
# ---------------------------------------------------------------------
$bmIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*PATIENTFIRSTNAME*") {
        $bmIndex = $i
        break
    }
}
$bmPara = $d.Paragraphs.Item($bmIndex)
$bmPara.Range.InsertParagraphBefore()
$bmPara.Range.InsertParagraphBefore()
$bmPara.Range.InsertParagraphBefore()

# re-locate after the inserts (paragraph identity/index shifted)
$bmIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*PATIENTFIRSTNAME*") {
        $bmIndex = $i
        break
    }
}

$nameTargetPara   = $d.Paragraphs.Item($bmIndex - 3)
$agencyTargetPara = $d.Paragraphs.Item($bmIndex - 2)
$genderTargetPara = $d.Paragraphs.Item($bmIndex - 1)
$bmPara           = $d.Paragraphs.Item($bmIndex)
$oldCombinedPara  = $d.Paragraphs.Item($bmIndex + 1)

# "Patient " (bold) + patient name
$r1 = $nameTargetPara.Range
$t1 = $d.Range($r1.Start, $r1.End - 1)
$t1.Text = "Patient PATIENTFIRSTNAME GOES HERE PATIENTLASTNAME GOES HERE"
$boldPatient = $d.Range($t1.Start, $t1.Start + 8)
$boldPatient.Bold = 1

# Agency ID (bare value only, no label/comma)
$r2 = $agencyTargetPara.Range
$t2 = $d.Range($r2.Start, $r2.End - 1)
$t2.Text = "PATIENTAGENCYID GOES HERE"

# Gender (bare value only, no comma)
$r3 = $genderTargetPara.Range
$t3 = $d.Range($r3.Start, $r3.End - 1)
$t3.Text = "PATIENTGENDER GOES HERE"

# bookmarked paragraph now just holds the age value
$rb = $bmPara.Range
$tb = $d.Range($rb.Start, $rb.End - 1)
$tb.Text = "PATIENTAGE GOES HERE"

# drop the now-redundant old combined paragraph's text
$oldCombinedPara.Range.Delete()

# ---------------------------------------------------------------------
# 3. Update the hard-coded date text
# ---------------------------------------------------------------------
$d.Content.Find.Execute("23 March 2013", $true, $false, $false, $false, $false, $true, 1, $false, "3 April 2013", 2) | Out-Null
